$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep a text representation (they are stored as inline strings
# in the original workbook), so force text number format before assigning
# values that look numeric (otherwise Excel auto-converts them to numbers).
$ws.Range("C2:F2").NumberFormat = "@"
$ws.Range("M2:T2").NumberFormat = "@"

$ws.Range("A2").Value = "2025-03-26T12:33"
$ws.Range("B2").Value = "SHIFT_1"
$ws.Range("C2").Value = "0.2"
$ws.Range("D2").Value = "0.2"
$ws.Range("E2").Value = "0.2"
$ws.Range("F2").Value = "0.3"

$ws.Range("M2").Value = "300"
$ws.Range("N2").Value = "200"
$ws.Range("O2").Value = "0.2"
$ws.Range("P2").Value = "1.23"
$ws.Range("Q2").Value = "1.6"
$ws.Range("R2").Value = "44"
$ws.Range("S2").Value = "2.63"
$ws.Range("T2").Value = "32.3"

$ws.Range("W2").Value = "Suriya"
